# Daily attendance processing - reorder "Recorded By" (column G) values so
# that "System" appears first in the comma-separated list instead of last.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value()

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }
    if ($val -notlike "*,*") { continue }

    $parts = $val -split ","
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts[$parts.Length - 1] -eq "System") {
        $rest = $parts[0..($parts.Length - 2)]
        $newVal = "System, " + ($rest -join ", ")
        $cell.Value = $newVal
    }
}
